$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.359579086303711
$ws.Range("B1").Value = 3.689178943634033
$ws.Range("C1").Value = 0.8204050660133362
$ws.Range("D1").Value = 0.4964723289012909
$ws.Range("E1").Value = 0.1406157165765762
